$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert new rows (bottom-up so row numbers below are not yet shifted) ---
$ws.Range("10:15").Insert()   # room for 6 new Sergipe years (2015-2020)
$ws.Range("6:11").Insert()    # room for 6 new Nordeste years (2015-2020)
$ws.Range("2:7").Insert()     # room for 6 new Brasil years (2015-2020)

# Inserting directly below the styled header row causes Excel to copy the
# headers style onto the freshly inserted blank rows; strip that back to Normal
# so the new data cells end up with the default (no-style) formatting.
$ws.Range("2:7").Style = "Normal"

# --- Step 2: new column F header, formatted like the other header cells ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Faltam dados para todos os Estados"

# --- Step 3: fill in the data rows ---
# row 2: Brasil 01/01/2015
$ws.Cells.Item(2, 1).Value = "Brasil"
$c = $ws.Cells.Item(2, 2)
$c.Value = "'01/01/2015"
$c.Style = "Normal"
$ws.Cells.Item(2, 3).Value = "Roubo de veículo"
$ws.Cells.Item(2, 4).Value = 84.8415561828753
$ws.Cells.Item(2, 6).Value = $true

# row 3: Brasil 01/01/2016
$ws.Cells.Item(3, 1).Value = "Brasil"
$c = $ws.Cells.Item(3, 2)
$c.Value = "'01/01/2016"
$c.Style = "Normal"
$ws.Cells.Item(3, 3).Value = "Roubo de veículo"
$ws.Cells.Item(3, 4).Value = 102.8390941398115
$ws.Cells.Item(3, 6).Value = $true

# row 4: Brasil 01/01/2017
$ws.Cells.Item(4, 1).Value = "Brasil"
$c = $ws.Cells.Item(4, 2)
$c.Value = "'01/01/2017"
$c.Style = "Normal"
$ws.Cells.Item(4, 3).Value = "Roubo de veículo"
$ws.Cells.Item(4, 4).Value = 108.1721960914802
$ws.Cells.Item(4, 6).Value = $true

# row 5: Brasil 01/01/2018
$ws.Cells.Item(5, 1).Value = "Brasil"
$c = $ws.Cells.Item(5, 2)
$c.Value = "'01/01/2018"
$c.Style = "Normal"
$ws.Cells.Item(5, 3).Value = "Roubo de veículo"
$ws.Cells.Item(5, 4).Value = 102.2104007004055
$ws.Cells.Item(5, 6).Value = $true

# row 6: Brasil 01/01/2019
$ws.Cells.Item(6, 1).Value = "Brasil"
$c = $ws.Cells.Item(6, 2)
$c.Value = "'01/01/2019"
$c.Style = "Normal"
$ws.Cells.Item(6, 3).Value = "Roubo de veículo"
$ws.Cells.Item(6, 4).Value = 83.51050151510037
$ws.Cells.Item(6, 6).Value = $true

# row 7: Brasil 01/01/2020
$ws.Cells.Item(7, 1).Value = "Brasil"
$c = $ws.Cells.Item(7, 2)
$c.Value = "'01/01/2020"
$c.Style = "Normal"
$ws.Cells.Item(7, 3).Value = "Roubo de veículo"
$ws.Cells.Item(7, 4).Value = 69.45399507865673
$ws.Cells.Item(7, 6).Value = $true

# row 8: Brasil 01/01/2021
$ws.Cells.Item(8, 6).Value = $false

# row 9: Brasil 01/01/2022
$ws.Cells.Item(9, 6).Value = $false

# row 10: Brasil 01/01/2023
$ws.Cells.Item(10, 6).Value = $false

# row 11: Brasil 01/01/2024
$ws.Cells.Item(11, 6).Value = $false

# row 12: Nordeste 01/01/2015
$ws.Cells.Item(12, 1).Value = "Nordeste"
$c = $ws.Cells.Item(12, 2)
$c.Value = "'01/01/2015"
$c.Style = "Normal"
$ws.Cells.Item(12, 3).Value = "Roubo de veículo"
$ws.Cells.Item(12, 4).Value = 94.18584801468253
$ws.Cells.Item(12, 6).Value = $true

# row 13: Nordeste 01/01/2016
$ws.Cells.Item(13, 1).Value = "Nordeste"
$c = $ws.Cells.Item(13, 2)
$c.Value = "'01/01/2016"
$c.Style = "Normal"
$ws.Cells.Item(13, 3).Value = "Roubo de veículo"
$ws.Cells.Item(13, 4).Value = 116.8466589954399
$ws.Cells.Item(13, 6).Value = $true

# row 14: Nordeste 01/01/2017
$ws.Cells.Item(14, 1).Value = "Nordeste"
$c = $ws.Cells.Item(14, 2)
$c.Value = "'01/01/2017"
$c.Style = "Normal"
$ws.Cells.Item(14, 3).Value = "Roubo de veículo"
$ws.Cells.Item(14, 4).Value = 121.7601711824096
$ws.Cells.Item(14, 6).Value = $true

# row 15: Nordeste 01/01/2018
$ws.Cells.Item(15, 1).Value = "Nordeste"
$c = $ws.Cells.Item(15, 2)
$c.Value = "'01/01/2018"
$c.Style = "Normal"
$ws.Cells.Item(15, 3).Value = "Roubo de veículo"
$ws.Cells.Item(15, 4).Value = 114.7233992580749
$ws.Cells.Item(15, 6).Value = $true

# row 16: Nordeste 01/01/2019
$ws.Cells.Item(16, 1).Value = "Nordeste"
$c = $ws.Cells.Item(16, 2)
$c.Value = "'01/01/2019"
$c.Style = "Normal"
$ws.Cells.Item(16, 3).Value = "Roubo de veículo"
$ws.Cells.Item(16, 4).Value = 92.05996159991663
$ws.Cells.Item(16, 6).Value = $true

# row 17: Nordeste 01/01/2020
$ws.Cells.Item(17, 1).Value = "Nordeste"
$c = $ws.Cells.Item(17, 2)
$c.Value = "'01/01/2020"
$c.Style = "Normal"
$ws.Cells.Item(17, 3).Value = "Roubo de veículo"
$ws.Cells.Item(17, 4).Value = 94.22450942337706
$ws.Cells.Item(17, 6).Value = $true

# row 18: Nordeste 01/01/2021
$ws.Cells.Item(18, 6).Value = $false

# row 19: Nordeste 01/01/2022
$ws.Cells.Item(19, 6).Value = $false

# row 20: Nordeste 01/01/2023
$ws.Cells.Item(20, 6).Value = $false

# row 21: Nordeste 01/01/2024
$ws.Cells.Item(21, 6).Value = $false

# row 22: Sergipe 01/01/2015
$ws.Cells.Item(22, 1).Value = "Sergipe"
$c = $ws.Cells.Item(22, 2)
$c.Value = "'01/01/2015"
$c.Style = "Normal"
$ws.Cells.Item(22, 3).Value = "Roubo de veículo"
$ws.Cells.Item(22, 4).Value = 88.60189014358107
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = $true

# row 23: Sergipe 01/01/2016
$ws.Cells.Item(23, 1).Value = "Sergipe"
$c = $ws.Cells.Item(23, 2)
$c.Value = "'01/01/2016"
$c.Style = "Normal"
$ws.Cells.Item(23, 3).Value = "Roubo de veículo"
$ws.Cells.Item(23, 4).Value = 126.4565524072786
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = $true

# row 24: Sergipe 01/01/2017
$ws.Cells.Item(24, 1).Value = "Sergipe"
$c = $ws.Cells.Item(24, 2)
$c.Value = "'01/01/2017"
$c.Style = "Normal"
$ws.Cells.Item(24, 3).Value = "Roubo de veículo"
$ws.Cells.Item(24, 4).Value = 132.239620851065
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = $true

# row 25: Sergipe 01/01/2018
$ws.Cells.Item(25, 1).Value = "Sergipe"
$c = $ws.Cells.Item(25, 2)
$c.Value = "'01/01/2018"
$c.Style = "Normal"
$ws.Cells.Item(25, 3).Value = "Roubo de veículo"
$ws.Cells.Item(25, 4).Value = 99.19642120380563
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = $true

# row 26: Sergipe 01/01/2019
$ws.Cells.Item(26, 1).Value = "Sergipe"
$c = $ws.Cells.Item(26, 2)
$c.Value = "'01/01/2019"
$c.Style = "Normal"
$ws.Cells.Item(26, 3).Value = "Roubo de veículo"
$ws.Cells.Item(26, 4).Value = 82.51765408007823
$ws.Cells.Item(26, 5).Value = 11
$ws.Cells.Item(26, 6).Value = $true

# row 27: Sergipe 01/01/2020
$ws.Cells.Item(27, 1).Value = "Sergipe"
$c = $ws.Cells.Item(27, 2)
$c.Value = "'01/01/2020"
$c.Style = "Normal"
$ws.Cells.Item(27, 3).Value = "Roubo de veículo"
$ws.Cells.Item(27, 4).Value = 86.45848785182783
$ws.Cells.Item(27, 5).Value = 10
$ws.Cells.Item(27, 6).Value = $true

# row 28: Sergipe 01/01/2021
$ws.Cells.Item(28, 6).Value = $false

# row 29: Sergipe 01/01/2022
$ws.Cells.Item(29, 6).Value = $false

# row 30: Sergipe 01/01/2023
$ws.Cells.Item(30, 6).Value = $false

# row 31: Sergipe 01/01/2024
$ws.Cells.Item(31, 6).Value = $false

